# "Added my story cards" - fill in Asignee/Status for the three story rows
# that were still missing an owner (Search a File / Sort Files by Name /
# Sort Files by Date), matching the new shared-string entries "Braydon"
# and "Not Complete".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = "Braydon"
$ws.Range("E10").Value = "Not Complete"

$ws.Range("D14").Value = "Braydon"
$ws.Range("E14").Value = "Not Complete"

$ws.Range("D15").Value = "Braydon"
$ws.Range("E15").Value = "Not Complete"

# Minor cosmetic tweaks that were also present in the author's save
$ws.Columns("E").ColumnWidth = 13.856026785714286
$null = $ws.Range("F24").Select()
